$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert the "new changes in ops (ordercreation & orderpage & order form)" commit:
# header row 1
$ws.Range("C1").Value = "Emp ID-Order Assigned"
$ws.Range("D1").Value = "Assignee_QA"

# data row 2
$ws.Range("C2").Value = "SIPL5316"
$ws.Range("D2").Value = "SIPL5688"
$ws.Range("L2").Value = "WIP"

# cosmetic view-state restored to match the original template
$ws.Columns.Item(3).ColumnWidth = 20.33203125
$ws.Columns.Item(8).ColumnWidth = 28.88671875
$ws.Range("E2").Select()
